$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.6226591760299626
$summary.Range("C2").Value = 0.5746864310148233
$summary.Range("D2").Value = 0.9438202247191011
$summary.Range("E2").Value = 0.7143869596031184
$summary.Range("F2").Value = 0.8363757052771325
$summary.Range("G2").Value = 0.921065579531876
$summary.Range("H2").Value = 0.7968410273674761
$summary.Range("I2").Value = 504
$summary.Range("J2").Value = 373
$summary.Range("K2").Value = 161
$summary.Range("L2").Value = 30

# --- Classification Report sheet ---
$report = $wb.Worksheets.Item("Classification Report")

# class 0 row
$report.Range("B2").Value = 0.8429319371727748
$report.Range("C2").Value = 0.301498127340824
$report.Range("D2").Value = 0.4441379310344827

# class 1 row
$report.Range("B3").Value = 0.5746864310148233
$report.Range("C3").Value = 0.9438202247191011
$report.Range("D3").Value = 0.7143869596031184

# accuracy row
$report.Range("B4").Value = 0.6226591760299626
$report.Range("C4").Value = 0.6226591760299626
$report.Range("D4").Value = 0.6226591760299626
$report.Range("E4").Value = 0.6226591760299626

# macro avg row
$report.Range("B5").Value = 0.7088091840937991
$report.Range("C5").Value = 0.6226591760299626
$report.Range("D5").Value = 0.5792624453188006

# weighted avg row
$report.Range("B6").Value = 0.7088091840937991
$report.Range("C6").Value = 0.6226591760299626
$report.Range("D6").Value = 0.5792624453188006

# --- Confusion Matrix sheet ---
$cm = $wb.Worksheets.Item("Confusion Matrix")
$cm.Range("B2").Value = 161
$cm.Range("C2").Value = 373
$cm.Range("B3").Value = 30
$cm.Range("C3").Value = 504
